$d = $word.ActiveDocument

# The default footer (footer2.xml) has a second paragraph containing two
# tab characters followed by a DATE field ("M/d/yy"). Remove the date
# field and the leading tabs entirely, leaving an empty Footer-styled
# paragraph (commit: "removed date from footer because we dont want to
# incriminate ourselves").

$footer = $d.Sections(1).Footers(1)
$footerRange = $footer.Range

# Remove the DATE field (this deletes the begin/instrText/separate/
# result/end run sequence as a unit).
for ($i = $footerRange.Fields.Count; $i -ge 1; $i--) {
    $fld = $footerRange.Fields($i)
    if ($fld.Code.Text -match "DATE") {
        $fld.Delete()
    }
}

# Remove the two leading tab characters left behind in that paragraph.
$footerRange.Find.Execute("^t^t", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 2)
